$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two FAQ rows for "Cấu hình Vlan Switch HPE1111" / "Cấu hình Vlan Switch UPE5523"
# (rows 7 and 8), causing the rows below to shift up.
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(7).Delete()

# Update the active selection to reflect where the user ended up after editing (first blank row below data).
$ws.Range("B15").Select()
